# Daily IST report: add the 2026-02-19 column (CSV/MD/XLSX regen) (#393)
#
# Inserts a new "2026-02-19" submission-count column right after the
# existing "2026-02-18" column (pushing total_files/unique_days one
# column to the right, from H/I to I/J), fills in that day's per-person
# file counts, and rolls those counts into the total_files /
# unique_days summary columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new date column -------------------------------
# This shifts the old H (total_files) -> I and old I (unique_days) -> J,
# and (conveniently) the new H1 inherits H1's old header style while I1/J1
# keep their old total_files/unique_days (grey-fill) header style.
$ws.Columns("H:H").Insert()

# Match the target column width (stored width 12, same as the other date
# columns D:G). The ColumnWidth COM property is offset from the stored
# OOXML column width by the sheet's default padding (~0.8333 chars for
# this Calibri 11 workbook), so back that out here.
$ws.Columns("H:H").ColumnWidth = 11.166666666666666

# --- 2. Write the new header text without Excel's date auto-conversion --
# Typing "2026-02-19" straight into a cell makes Excel helpfully (and
# unhelpfully) reinterpret it as a real date. Stage it as Text in a
# scratch cell, then paste-special *values only* into H1 so the literal
# string lands without dragging the scratch cell's Text number format
# (or any other formatting) along with it.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "2026-02-19"
$ws.Range("Z1").Copy()
$ws.Range("H1").PasteSpecial(-4163)
$excel.CutCopyMode = 0
# Delete (not just clear) the scratch column so the sheet's used range /
# <dimension> doesn't balloon out to column Z.
$ws.Columns("Z:Z").Delete()

# --- 3. Per-row file counts submitted on 2026-02-19 ----------------------
$counts = @(1,0,1,1,0,1,1,1,1,1,1,1,1,0,0,1,1,1,0,0,0,0,0,1,0,0,0,0,0,1,1,1,1,1,1,0,0,0,1,1,1,0,0,0,1,1,0,1,1,0,1,1,0,1,0,1,0,1,1,0,0,0,1,0,0,0,0,0,1,0,0,0,1,0,0,0,0,0,1,1,1,0,0,0,0,0,1,0,0,0,0,0,1,1,0,0,0,0,0,1,0,1,0,1,1,0,0,0)

for ($i = 0; $i -lt $counts.Length; $i++) {
    $row = $i + 2
    $newDayCount = $counts[$i]

    # New day's column (H)
    $ws.Cells.Item($row, 8).Value = $newDayCount

    # total_files (now column I) picks up this day's files too
    $oldTotalFiles = $ws.Cells.Item($row, 9).Value()
    $ws.Cells.Item($row, 9).Value = $oldTotalFiles + $newDayCount

    # unique_days (now column J) gains one more active day if any files
    # were submitted on 2026-02-19
    $oldUniqueDays = $ws.Cells.Item($row, 10).Value()
    $dayIncrement = 0
    if ($newDayCount -gt 0) { $dayIncrement = 1 }
    $ws.Cells.Item($row, 10).Value = $oldUniqueDays + $dayIncrement
}
